$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 new columns before column U (old V/W shift right to Z/AA).
$ws.Columns("U:X").Insert()

# New header labels for the inserted columns (row 3 of the header block).
$ws.Range("U3").Value = "Running time (seg)"
$ws.Range("V3").Value = "Data volume (GB)"
$ws.Range("W3").Value = "Sample size"
$ws.Range("X3").Value = "Ranking"

# Extend the merged header cell R1:T2 to cover the newly added columns.
$ws.Range("R1:X2").Merge()
